# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# Both sheets list the same events; "全部类型" has one extra row (row 31,
# an event not shown on "展览"), so row numbers diverge after that point.

$wb = $excel.ActiveWorkbook

# row -> new F value, for the "展览" sheet (rows 2-34)
$exhibitionUpdates = @{
    2  = 112
    3  = 304
    5  = 614
    6  = 68
    7  = 2106
    10 = 4641
    11 = 3
    12 = 35
    14 = 213
    15 = 18
    16 = 148
    19 = 93
    20 = 3526
    21 = 87
    22 = 565
    24 = 21
    25 = 90
    27 = 12
    29 = 72
    30 = 213
    32 = 738
    33 = 2168
    34 = 400
}

# row -> new F value, for the "全部类型" sheet (rows 2-35)
$allTypesUpdates = @{
    2  = 112
    3  = 304
    5  = 614
    6  = 68
    7  = 2106
    10 = 4641
    11 = 3
    12 = 35
    14 = 213
    15 = 18
    16 = 148
    19 = 93
    20 = 3526
    21 = 87
    22 = 565
    24 = 21
    25 = 90
    27 = 12
    29 = 72
    30 = 213
    33 = 738
    34 = 2168
    35 = 400
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
